$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 270.1111   # was 296.46667
$ws.Range("I33").Value = 255.9   # was 285.29166
$ws.Range("K33").Value = 255.9   # was 285.29166
$ws.Range("M33").Value = -26.90000000000001   # was -56.29165999999998
$ws.Range("I40").Value = 1000   # was 0
$ws.Range("J40").Value = 1666.6666   # was 1500
$ws.Range("K40").Value = 1000   # was 0
$ws.Range("L40").Value = 1666.6666   # was 1500
$ws.Range("M40").Value = -825   # was None
$ws.Range("N40").Value = -2016.6666   # was -1850
$ws.Range("H93").Value = 0   # was 40000
$ws.Range("J93").Value = 0   # was 40000
$ws.Range("L93").Value = 0   # was 40000
$ws.Range("N93").Value = $null   # was -44992
$ws.Range("H116").Value = 1605.5   # was 1842.3636
$ws.Range("J116").Value = 1680.8572   # was 2253.2
$ws.Range("L116").Value = 1680.8572   # was 2253.2
$ws.Range("N116").Value = -8564.8572   # was -9137.200000000001
$ws.Range("H132").Value = 639293.75   # was 757302.2
$ws.Range("I132").Value = 2748.0308   # was 3388.5283
$ws.Range("J132").Value = 4087249.8   # was 4087087.5
$ws.Range("K132").Value = 8244.0924   # was 10165.5849
$ws.Range("L132").Value = 12261749.4   # was 12261262.5
$ws.Range("M132").Value = -5714.0924   # was -7635.5849
$ws.Range("N132").Value = -12266809.4   # was -12266322.5
$ws.Range("H138").Value = 4259154   # was 2781152.8
$ws.Range("I138").Value = 3640.6667   # was 2946.353
$ws.Range("J138").Value = 5718187   # was 3639871
$ws.Range("K138").Value = 10922.0001   # was 8839.059000000001
$ws.Range("L138").Value = 17154561   # was 10919613
$ws.Range("M138").Value = -5782.000100000001   # was -3699.059000000001
$ws.Range("N138").Value = -17164841   # was -10929893
$ws.Range("H141").Value = 2802.6775   # was 1605.7097
$ws.Range("I141").Value = 2616.6553   # was 1406.138
$ws.Range("J141").Value = 5500   # was 4499.5
$ws.Range("K141").Value = 7849.965899999999   # was 4218.414
$ws.Range("L141").Value = 16500   # was 13498.5
$ws.Range("M141").Value = -2669.965899999999   # was 961.5860000000002
$ws.Range("N141").Value = -26860   # was -23858.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 1000000   # was 334000
$ws.Range("J6").Value = 1000000   # was 334000
$ws.Range("L6").Value = 1000000   # was 334000
$ws.Range("N6").Value = -1000346   # was -334346
$ws.Range("H24").Value = 36000   # was 36588.75
$ws.Range("J24").Value = 36000   # was 36588.75
$ws.Range("L24").Value = 36000   # was 36588.75
$ws.Range("N24").Value = -36748   # was -37336.75
$ws.Range("H32").Value = 18485.3   # was 3781551.2
$ws.Range("I32").Value = 12754.189   # was 32416.41
$ws.Range("J32").Value = 34796.92   # was 10143719
$ws.Range("K32").Value = 12754.189   # was 32416.41
$ws.Range("L32").Value = 34796.92   # was 10143719
$ws.Range("M32").Value = -12467.189   # was -32129.41
$ws.Range("N32").Value = -35370.92   # was -10144293
$ws.Range("H100").Value = 36000   # was 36588.75
$ws.Range("J100").Value = 36000   # was 36588.75
$ws.Range("L100").Value = 36000   # was 36588.75
$ws.Range("N100").Value = -38164   # was -38752.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1705.2858   # was 2266.3333
$ws.Range("I107").Value = 1776.6923   # was 2266.3333
$ws.Range("J107").Value = 777   # was 0
$ws.Range("K107").Value = 1776.6923   # was 2266.3333
$ws.Range("L107").Value = 777   # was 0
$ws.Range("M107").Value = 143.3077000000001   # was -346.3332999999998
$ws.Range("N107").Value = -4617   # was None

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 22224758   # was 83336584
$ws.Range("I58").Value = 29414110   # was 166669630
$ws.Range("J58").Value = 3125.4546   # was 3531.6667
$ws.Range("K58").Value = 29414110   # was 166669630
$ws.Range("L58").Value = 3125.4546   # was 3531.6667
$ws.Range("M58").Value = -29413907   # was -166669427
$ws.Range("N58").Value = -3531.4546   # was -3937.6667
$ws.Range("H132").Value = 53007.95   # was 53006.85
$ws.Range("I132").Value = 2797.2144   # was 2795.6428
$ws.Range("K132").Value = 8391.643199999999   # was 8386.928400000001
$ws.Range("M132").Value = -5861.643199999999   # was -5856.928400000001
$ws.Range("H134").Value = 270653   # was 99652.27
$ws.Range("I134").Value = 1112   # was 1335
$ws.Range("J134").Value = 360500   # was 181583.33
$ws.Range("K134").Value = 3336   # was 4005
$ws.Range("L134").Value = 1081500   # was 544749.99
$ws.Range("M134").Value = -801   # was -1470
$ws.Range("N134").Value = -1086570   # was -549819.99
$ws.Range("H136").Value = 22224758   # was 83336584
$ws.Range("I136").Value = 29414110   # was 166669630
$ws.Range("J136").Value = 3125.4546   # was 3531.6667
$ws.Range("K136").Value = 88242330   # was 500008890
$ws.Range("L136").Value = 9376.363799999999   # was 10595.0001
$ws.Range("M136").Value = -88239780   # was -500006340
$ws.Range("N136").Value = -14476.3638   # was -15695.0001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 828.7692   # was 730.7222
$ws.Range("I113").Value = 594   # was 592.8570999999999
$ws.Range("J113").Value = 975.5   # was 818.4545000000001
$ws.Range("K113").Value = 1782   # was 1778.5713
$ws.Range("L113").Value = 2926.5   # was 2455.3635
$ws.Range("M113").Value = 388   # was 391.4287000000002
$ws.Range("N113").Value = -7266.5   # was -6795.3635
$ws.Range("H137").Value = 3558.1177   # was 2889.4119
$ws.Range("I137").Value = 1363.3334   # was 1024
$ws.Range("J137").Value = 4028.4285   # was 3666.6667
$ws.Range("K137").Value = 4090.0002   # was 3072
$ws.Range("L137").Value = 12085.2855   # was 11000.0001
$ws.Range("M137").Value = 1009.9998   # was 2028
$ws.Range("N137").Value = -22285.2855   # was -21200.0001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1858.9656   # was 2099.8
$ws.Range("I113").Value = 1131.2941   # was 1279.6428
$ws.Range("J113").Value = 2889.8333   # was 3143.6365
$ws.Range("K113").Value = 1131.2941   # was 1279.6428
$ws.Range("L113").Value = 2889.8333   # was 3143.6365
$ws.Range("M113").Value = 1038.7059   # was 890.3571999999999
$ws.Range("N113").Value = -7229.8333   # was -7483.636500000001
$ws.Range("H132").Value = 51411.926   # was 62947.547
$ws.Range("I132").Value = 32421.182   # was 37980.605
$ws.Range("J132").Value = 129748.75   # was 202762.4
$ws.Range("K132").Value = 97263.546   # was 113941.815
$ws.Range("L132").Value = 389246.25   # was 608287.2
$ws.Range("M132").Value = -94733.546   # was -111411.815
$ws.Range("N132").Value = -394306.25   # was -613347.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 842.7895   # was 833.55554
$ws.Range("I22").Value = 751.5   # was 714.2857
$ws.Range("J22").Value = 1098.4   # was 1251
$ws.Range("K22").Value = 751.5   # was 714.2857
$ws.Range("L22").Value = 1098.4   # was 1251
$ws.Range("M22").Value = -456.5   # was -419.2857
$ws.Range("N22").Value = -1688.4   # was -1841
$ws.Range("H27").Value = 842.7895   # was 833.55554
$ws.Range("I27").Value = 751.5   # was 714.2857
$ws.Range("J27").Value = 1098.4   # was 1251
$ws.Range("K27").Value = 751.5   # was 714.2857
$ws.Range("L27").Value = 1098.4   # was 1251
$ws.Range("M27").Value = -644.5   # was -607.2857
$ws.Range("N27").Value = -1312.4   # was -1465
$ws.Range("H93").Value = 2161.6   # was 2200.8
$ws.Range("J93").Value = 2202   # was 2251
$ws.Range("L93").Value = 2202   # was 2251
$ws.Range("N93").Value = -4698   # was -4747
$ws.Range("H136").Value = 51095.49   # was 51148.7
$ws.Range("I136").Value = 31383.81   # was 31445.648
$ws.Range("K136").Value = 94151.43000000001   # was 94336.944
$ws.Range("M136").Value = -91601.43000000001   # was -91786.944

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 673.2105   # was 640.5238000000001
$ws.Range("I107").Value = 526.8   # was 503.64706
$ws.Range("K107").Value = 1580.4   # was 1510.94118
$ws.Range("M107").Value = 339.6000000000001   # was 409.05882
$ws.Range("H113").Value = 824.3214   # was 921.125
$ws.Range("I113").Value = 989   # was 1188.4166
$ws.Range("J113").Value = 634.3077   # was 653.8333
$ws.Range("K113").Value = 2967   # was 3565.2498
$ws.Range("L113").Value = 1902.9231   # was 1961.4999
$ws.Range("M113").Value = -797   # was -1395.2498
$ws.Range("N113").Value = -6242.9231   # was -6301.4999
$ws.Range("H122").Value = 2427.647   # was 2505
$ws.Range("I122").Value = 1624.2858   # was 1667.1428
$ws.Range("J122").Value = 2990   # was 3238.125
$ws.Range("K122").Value = 4872.857400000001   # was 5001.428400000001
$ws.Range("L122").Value = 8970   # was 9714.375
$ws.Range("M122").Value = -2422.857400000001   # was -2551.428400000001
$ws.Range("N122").Value = -13870   # was -14614.375
$ws.Range("H132").Value = 44985.78   # was 48087.023
$ws.Range("I132").Value = 38109.15   # was 39567.04
$ws.Range("J132").Value = 54757.844   # was 61117.59
$ws.Range("K132").Value = 114327.45   # was 118701.12
$ws.Range("L132").Value = 164273.532   # was 183352.77
$ws.Range("M132").Value = -111797.45   # was -116171.12
$ws.Range("N132").Value = -169333.532   # was -188412.77
